# edit.ps1 - applies the "Faza 1" revision to "Agenti de vanzari.docx"
#
# Strategy notes (from experimentation against this COM-interop runtime):
#  - Range.Find.Execute(...) narrows a Range to the matched text, which we
#    then edit in place.
#  - Range.Text = "..." replaces the matched text but keeps the original
#    run's formatting (rPr) - perfect for simple retext / "collapse multiple
#    runs into one run" edits.
#  - To split one run into several runs that have genuinely *different*
#    formatting (e.g. add <w:u w:val="single"/> on one word), we just set
#    the desired Font property on the relevant sub-Range; the engine
#    automatically creates a new run boundary there.
#  - To split a run into several runs that end up with *identical*
#    formatting (pure text split, no property change), setting a property
#    and reverting it leaves a stray explicit value in rPr. Instead we
#    round-trip Range.FormattedText (read then write back the same value)
#    on the sub-range, which forces a clean run boundary with no residue.
#  - New list paragraphs are created with Range.InsertAfter(Chr(13) + text);
#    the runtime automatically carries over the source paragraph's pPr
#    (style + numbering) to the freshly inserted paragraph.

$d = $word.ActiveDocument

function Split-Run($rangeStart, $rangeEnd) {
    # Forces a run boundary at [rangeStart,rangeEnd) without altering
    # formatting, by round-tripping FormattedText.
    $sub = $d.Range($rangeStart, $rangeEnd)
    $ft = $sub.FormattedText
    $sub.FormattedText = $ft
}

# ---------------------------------------------------------------------
# 1) Title: "Agenti de vanzari" -> "2. " + "Agenti de vanzari" (2 runs)
# ---------------------------------------------------------------------
$full = $d.Content
$find = $full.Find
$find.ClearFormatting()
$target = "Agenti de vanzari"
$found = $find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $s = $full.Start
    $full.Text = "2. " + $target
    $splitStart = $s + 3
    $splitEnd = $s + 3 + $target.Length
    Split-Run $splitStart $splitEnd
}

# ---------------------------------------------------------------------
# 2) "...impreuna cu preturile aferente si cantitatile existente pe stoc."
#    -> underline "preturile" and "cantitatile"
# ---------------------------------------------------------------------
$full = $d.Content
$find = $full.Find
$find.ClearFormatting()
$target = 'agentul vizualizeaza lista tuturor produselor vândute de firma, împreuna cu preturile aferente și cantitatile existente pe stoc. '
$found = $find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $s = $full.Start
    $full.Text = $target

    $w1 = "preturile"
    $off1 = $target.IndexOf($w1)
    $sub1 = $d.Range($s + $off1, $s + $off1 + $w1.Length)
    $sub1.Font.Underline = 1

    $w2 = "cantitatile"
    $off2 = $target.IndexOf($w2)
    $sub2 = $d.Range($s + $off2, $s + $off2 + $w2.Length)
    $sub2.Font.Underline = 1
}

# ---------------------------------------------------------------------
# 3) "...Dupa orice comanda valida, toti agentii logati  in  aplicatie..."
#    -> underline "comanda valida" and "logati"
# ---------------------------------------------------------------------
$full = $d.Content
$find = $full.Find
$find.ClearFormatting()
$target = 'agentul poate comanda o cantitate dintr-un produs. După orice comanda valida, toti agentii logati  în  aplicație  vor  vedea  lista  actualizata  a  stocurilor  (este  posibil  ca  declansarea  înregistrarii unei comenzi din partea unui agent sa conduca la un mesaj informativ "cantitate insuficienta în stoc").'
$found = $find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $s = $full.Start
    $full.Text = $target

    $w1 = "comanda valida"
    $off1 = $target.IndexOf($w1)
    $sub1 = $d.Range($s + $off1, $s + $off1 + $w1.Length)
    $sub1.Font.Underline = 1

    $w2 = "logati"
    $off2 = $target.IndexOf($w2)
    $sub2 = $d.Range($s + $off2, $s + $off2 + $w2.Length)
    $sub2.Font.Underline = 1
}

# ---------------------------------------------------------------------
# 4) "Logare client/ Creare cont nou client" -> "Logare " + "agent" (2 runs)
# ---------------------------------------------------------------------
$full = $d.Content
$find = $full.Find
$find.ClearFormatting()
$target = "Logare client/ Creare cont nou client"
$found = $find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $s = $full.Start
    $newText = "Logare agent"
    $full.Text = $newText
    $splitStart = $s + 7
    $splitEnd = $s + $newText.Length
    Split-Run $splitStart $splitEnd
}

# ---------------------------------------------------------------------
# 5) "Filtrare" + "/sortare" + " " + "tip/pret produse" (4 runs)
#    -> "Vizualizare produse/comenzi" (1 run)
# ---------------------------------------------------------------------
$full = $d.Content
$find = $full.Find
$find.ClearFormatting()
$target = "Filtrare/sortare tip/pret produse"
$found = $find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $full.Text = "Vizualizare produse/comenzi"
}

# ---------------------------------------------------------------------
# 6) "Oferte valabile" -> "Cauta produse"
# ---------------------------------------------------------------------
$full = $d.Content
$find = $full.Find
$find.ClearFormatting()
$target = "Oferte valabile"
$found = $find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $full.Text = "Cauta produse"
}

# ---------------------------------------------------------------------
# 7) "Plasare/Anulare comanda" -> "Adauga in cos"
#    then add two new list paragraphs:
#      "Vizualizeaza cos"
#      "Plaseaza/Anuleaza comanda"
# ---------------------------------------------------------------------
$full = $d.Content
$find = $full.Find
$find.ClearFormatting()
$target = "Plasare/Anulare comanda"
$found = $find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $full.Text = "Adauga in cos"
}

$full = $d.Content
$find = $full.Find
$find.ClearFormatting()
$target = "Adauga in cos"
$found = $find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $full.InsertAfter([char]13 + "Vizualizeaza cos" + [char]13 + "Plaseaza/Anuleaza comanda")
}
